$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Implement real electrical net content for the simplest set of electric parameters*") {
        $p.Range.Delete()
        break
    }
}
